$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 36: add "Pending" status in column B (bold, gold/orange font - matches other Pending cells)
$ws.Range("B36").Value = "Pending"
$ws.Range("B36").Font.Bold = $true
$ws.Range("B36").Font.Color = 49407

# Row 37: add "Pending" status in column B
$ws.Range("B37").Value = "Pending"
$ws.Range("B37").Font.Bold = $true
$ws.Range("B37").Font.Color = 49407

# New row 38: Task 38, with Pending status in column B
$ws.Range("A38").Value = "Task 38: Manage project page"
$ws.Range("B38").Value = "Pending"
$ws.Range("B38").Font.Bold = $true
$ws.Range("B38").Font.Color = 49407

# Update the view's top-left cell and selection to match author's final view
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B45").Select()
